# Trade #6 closed at 2026-02-17 04:06:21 on the MarketMaking strategy.
# Update the rollup sheets (Summary, Strategy Status) and append the new
# trade row to the detail sheets (All Trades, MarketMaking).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.97   # Current Capital
$summary.Range("B4").Value = -0.03     # Total P&L $
$summary.Range("B5").Value = -0.1      # Total P&L %
$summary.Range("B6").Value = 6         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking is row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.97   # Capital
$status.Range("D4").Value = 6       # Trades
$status.Range("E4").Value = -0.03   # P&L $
$status.Range("F4").Value = -0.03   # P&L %
$status.Range("G4").Value = 33.33   # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade (#6) to both the "All Trades" log and the
# per-strategy "MarketMaking" sheet - they mirror each other.
# ---------------------------------------------------------------------
$newRow = 7

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$newRow").Value = 6
    # Leading apostrophe forces plain text so "2026-02-17" isn't
    # auto-converted into a date serial number (matches the other rows,
    # which store the date as literal text).
    $ws.Range("B$newRow").Value = "'2026-02-17"
    $ws.Range("C$newRow").Value = "04:06:21"
    $ws.Range("D$newRow").Value = "MarketMaking"
    $ws.Range("E$newRow").Value = "DOWN"
    $ws.Range("F$newRow").Value = 0.78
    $ws.Range("G$newRow").Value = 0.76
    $ws.Range("H$newRow").Value = "CLOSED"
    $ws.Range("I$newRow").Value = -2.5641
    $ws.Range("J$newRow").Value = -0.02
    $ws.Range("K$newRow").Value = 99.97
    $ws.Range("L$newRow").Value = 0
    $ws.Range("M$newRow").Value = 0
    $ws.Range("N$newRow").Value = 0.6
    $ws.Range("O$newRow").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$newRow").Value = "early_exit"
    $ws.Range("Q$newRow").Value = 0.12
}
